$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 216; this shifts the existing rows 216-275 down to 217-276
# and copies formatting (including the date number format on column D) from the row above.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new record's data.
$ws.Range("A216").Value = 5
$ws.Range("B216").Value = "Macroferia Regional de Talca"
$ws.Range("C216").Value = "Maule"
$ws.Range("D216").Value = 44785
$ws.Range("E216").Value = 7
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100108
$ws.Range("H216").Value = "Tropicales y subtropicales"
$ws.Range("I216").Value = 100108005
$ws.Range("J216").Value = "Piña"
$ws.Range("K216").Value = "Caramelo"
$ws.Range("L216").Value = "Segunda"
$ws.Range("M216").Value = 200
$ws.Range("N216").Value = 19000
$ws.Range("O216").Value = 19000
$ws.Range("P216").Value = 19000
$ws.Range("Q216").Value = "$/caja 14 unidades"
$ws.Range("R216").Value = "Ecuador"
$ws.Range("S216").Value = 1357
$ws.Range("T216").Value = 14
